$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column I to fit the new, longer note text
$ws.Columns.Item(9).ColumnWidth = 33.33

# Fill in row 18: new binary size measurement + release note
# (match the date-formatted style already used by the rows above)
$ws.Range("A17").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("A18").Value = 45472.625694444447
$ws.Range("B18").Value = 1438720
$ws.Range("C18").Value = 259584
$ws.Range("D18").Value = 337408
$ws.Range("I18").Value = "v0.1.0.656 (v0.1.0, prerelease 1) RC1"

# Move the active selection to L24, matching the final saved view state
$ws.Range("L24").Select()
